$wb = $excel.ActiveWorkbook

# --- 1) Generic sheet: NrBuckets goes from 4 to 5 ---
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 5

# --- 2) ForecastedAverageDemand: add bucket row 6 (copy formatting from row 5, then set values) ---
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvg.Range("A5:N5").Copy()
$wsAvg.Range("A6").PasteSpecial(-4122)
$wsAvg.Cells.Item(6,1).Value = 4
$wsAvg.Cells.Item(6,2).Value = 0
$wsAvg.Cells.Item(6,3).Value = 0
$wsAvg.Cells.Item(6,4).Value = 0
$wsAvg.Cells.Item(6,5).Value = 0
$wsAvg.Cells.Item(6,6).Value = 0
$wsAvg.Cells.Item(6,7).Value = 4000
$wsAvg.Cells.Item(6,8).Value = 2000
$wsAvg.Cells.Item(6,9).Value = 8000
$wsAvg.Cells.Item(6,10).Value = 700
$wsAvg.Cells.Item(6,11).Value = 0
$wsAvg.Cells.Item(6,12).Value = 0
$wsAvg.Cells.Item(6,13).Value = 0
$wsAvg.Cells.Item(6,14).Value = 0

# --- 3) ForcastedStandardDeviation: add bucket row 6 (copy formatting from row 5, then set values) ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStd.Range("A5:N5").Copy()
$wsStd.Range("A6").PasteSpecial(-4122)
$wsStd.Cells.Item(6,1).Value = 4
$wsStd.Cells.Item(6,2).Value = 0
$wsStd.Cells.Item(6,3).Value = 0
$wsStd.Cells.Item(6,4).Value = 0
$wsStd.Cells.Item(6,5).Value = 0
$wsStd.Cells.Item(6,6).Value = 0
$wsStd.Cells.Item(6,7).Value = 7500
$wsStd.Cells.Item(6,8).Value = 3500
$wsStd.Cells.Item(6,9).Value = 9000
$wsStd.Cells.Item(6,10).Value = 1400
$wsStd.Cells.Item(6,11).Value = 0
$wsStd.Cells.Item(6,12).Value = 0
$wsStd.Cells.Item(6,13).Value = 0
$wsStd.Cells.Item(6,14).Value = 0
